$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '67.097.70'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  -1.03%  '

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.297.19'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  +0.83%  '

$ws.Range('E4').Value = '  -0.04%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '186.08'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +1.11%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '577.19'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.48%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.999'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -0.04%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.600'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -0.27%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.130'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -0.23%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '6.67'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +1.35%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.411'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +0.70%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '3.860.27'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +0.55%  '

$ws.Range('E13').Value = '  -0.53%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '27.48'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +0.35%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '67.354.34'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -0.66%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.0000167'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -0.60%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '3.272.82'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +0.42%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '444.08'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +10.41%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '5.71'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +0.15%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '13.54'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +0.97%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '7.75'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +2.70%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '74.43'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +4.61%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.999'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -0.07%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '0.514'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +1.40%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '3.425.65'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +0.46%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.0000119'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +1.49%  '

$ws.Range('E27').Value = '  -0.49%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '9.11'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -4.00%  '

$ws.Range('E29').Value = '  -0.06%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.97'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +1.28%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '22.83'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +0.53%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '5.35'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -2.23%  '

$ws.Range('E33').Value = '  -0.01%  '

$ws.Range('B34').Value = 'Fetch.AI'
$ws.Range('C34').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.24'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -0.40%  '

$ws.Range('B35').Value = 'Aptos'
$ws.Range('C35').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '6.81'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -1.50%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.53'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +5.38%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '162.34'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -1.28%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '27.42'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +1.15%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '1.85'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -2.02%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '4.48'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -0.34%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.783'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -2.45%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '2.729.73'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +1.96%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '6.29'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -0.74%  '

$ws.Range('B44').Value = 'Hedera'
$ws.Range('C44').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.0672'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -0.83%  '

$ws.Range('B45').Value = 'OKB'
$ws.Range('C45').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '40.16'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -1.57%  '

$ws.Range('B46').Value = 'InjectiveProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '24.80'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +0.78%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '2.41'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -1.13%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '328.96'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -2.08%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.0274'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -0.05%  '

$ws.Range('B50').Value = 'ONDO'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.991'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +2.38%  '

$ws.Range('B51').Value = 'Cosmos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '6.22'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -0.98%  '
